$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.354.75"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +3.85%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.718.09"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +3.29%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'0.9988"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.01%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'239.24"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.37%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'0.9995"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -0.02%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4708"
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').Value = "'0.2635"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.91%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.06219"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.23%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'1.710.91"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +2.79%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'0.07070"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.14%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'15.21"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +3.35%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.5900"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -0.42%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'4.415"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +0.77%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'76.33"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.82%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').Value = "'0.9997"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +0.03%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'0.9999"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.07%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'26.349.47"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +3.79%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = "'0.000006814"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.36%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.97%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'1.931.50"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +3.23%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'4.542"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +2.39%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'8.811"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.89%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'5.349"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.32%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'135.41"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.25%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'15.19"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.73%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'1.405"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +0.54%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'1.763"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +4.11%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'106.83"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +2.27%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'4.041"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +1.65%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'3.686"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +2.03%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'0.07706"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.68%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.04425"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +0.37%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'2.612"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.43%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'0.6217"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +2.42%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'0.9709"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +3.09%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.9340"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +9.47%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'114.40"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +16.01%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('D39').Value = "'2.404"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  -8.41%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = "'PaxDollar"
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = "'1.000"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +0.01%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('B41').Value = "'RenderToken"
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = "'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = "'1.906"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +4.45%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.01464"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -2.61%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'5.278"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +12.21%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'0.3803"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.00%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.1147"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +2.77%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'6.238"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.48%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'0.05287"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +0.71%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'30.50"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +3.43%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'7.679"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +5.20%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'1.221"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.47%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.3374"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.84%  "
$ws.Range('E51').Style = 'Normal'
